# 自动更新Excel文件 - daily roll-forward of "剩余" (remaining days) counter.
#
# Business rule (one row per shop, columns: D=总天 total days,
# E=剩余 remaining days, F=开始时间 start date as YYYYMMDD number):
#   - For every data row (skip the header in row 1):
#       * if 剩余 (E) has hit 1 (about to run out), the shop is restocked:
#           E is reset back to the total (D) and F (start date) is bumped
#           to the new restock date 2025-12-22 (20251222).
#       * otherwise one more day has elapsed, so E is simply decremented
#         by 1, leaving F untouched.
#   - Rows whose start date isn't a well-formed YYYYMMDD number (data
#     corruption) are left completely alone, matching upstream's skip
#     behaviour for unparsable dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$restockDate = 20251222

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $total = $ws.Cells.Item($r, 4).Value2
    $remaining = $ws.Cells.Item($r, 5).Value2
    $startDate = $ws.Cells.Item($r, 6).Value2

    if ($total -eq $null -or $remaining -eq $null) {
        continue
    }

    # Skip rows whose start-date cell isn't a clean 8-digit YYYYMMDD
    # number (corrupted source data) - they are left untouched.
    $dateText = [string]$startDate
    if ($dateText.Length -ne 8) {
        continue
    }

    if ($remaining -eq 1) {
        $ws.Cells.Item($r, 5).Value2 = $total
        $ws.Cells.Item($r, 6).Value2 = $restockDate
    } else {
        $ws.Cells.Item($r, 5).Value2 = $remaining - 1
    }
}
